# Insert a new weekly price record for "Feria Lagunitas de Puerto Montt - Cilantro"
# at row 181, pushing the existing rows 181-191 down to 182-192.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 181 (shifts 181:191 -> 182:192)
$ws.Rows(181).Insert()

# Populate the newly inserted row 181 with the new weekly data point
$ws.Range("A181").Value = 4
$ws.Range("B181").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C181").Value = "Los Lagos"
$ws.Range("D181").Value = 44516
$ws.Range("D181").NumberFormat = $ws.Range("D182").NumberFormat
$ws.Range("E181").Value = 10
$ws.Range("F181").Value = 100112040
$ws.Range("G181").Value = "Cilantro"
$ws.Range("H181").Value = "Sin especificar"
$ws.Range("I181").Value = "Primera"
$ws.Range("J181").Value = 250
$ws.Range("K181").Value = 10000
$ws.Range("L181").Value = 10000
$ws.Range("M181").Value = 10000
$ws.Range("N181").Value = "$/caja 36 atados"
$ws.Range("O181").Value = "Región Metropolitana"
$ws.Range("P181").Value = 278
$ws.Range("Q181").Value = 36
$ws.Range("R181").Value = "Hortaliza"
